# Update column G ("K") values for rows 2-13 to reflect the regenerated
# save_data using K instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 8
    4  = 0
    5  = 5
    6  = 4
    7  = 3
    8  = 6
    9  = 3
    10 = 1
    11 = 4
    12 = 6
    13 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
